$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column F (dSF) values for rows 2-10 as per repulled data
$ws.Range("F2").Value = 9
$ws.Range("F3").Value = -8
$ws.Range("F4").Value = -2
$ws.Range("F5").Value = 1
$ws.Range("F7").Value = 1
$ws.Range("F8").Value = -2
$ws.Range("F9").Value = -5
$ws.Range("F10").Value = -1
